$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "test@email.com"
$ws.Range("B4").Value = "'1"
$ws.Range("C4").Value = "'2024-11-20"
$ws.Range("D4").Value = "'21:52:17"

$ws.Range("A5").Value = "myemail45@gmail.com"
$ws.Range("B5").Value = "'1"
$ws.Range("C5").Value = "'2024-11-20"
$ws.Range("D5").Value = "'22:41:50"

$ws.Range("B4:D5").Style = "Normal"
